$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 370; everything from old row 370 downward
# shifts down by one (old 370 -> new 371, ..., old 492 -> new 493).
$ws.Rows.Item(370).Insert()

# Populate the newly inserted row 370 with the new weekly record
# (same dimension/category columns as the old row 370, new date + prices).
$ws.Cells.Item(370, 1).Value = 4
$ws.Cells.Item(370, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(370, 3).Value = "Los Lagos"
$ws.Cells.Item(370, 4).Value = 44988
$ws.Cells.Item(370, 5).Value = 10
$ws.Cells.Item(370, 6).Value = 100112008
$ws.Cells.Item(370, 7).Value = "Coliflor"
$ws.Cells.Item(370, 8).Value = "Sin especificar"
$ws.Cells.Item(370, 9).Value = "Primera"
$ws.Cells.Item(370, 10).Value = 750
$ws.Cells.Item(370, 11).Value = 1700
$ws.Cells.Item(370, 12).Value = 1700
$ws.Cells.Item(370, 13).Value = 1700
$ws.Cells.Item(370, 14).Value = "`$/unidad"
$ws.Cells.Item(370, 15).Value = "Región Metropolitana"
$ws.Cells.Item(370, 16).Value = 1700
$ws.Cells.Item(370, 17).Value = 1
$ws.Cells.Item(370, 18).Value = "Hortaliza"
